$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 8000
$ws.Cells.Item(13, 10).Value = 8000
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(43, 8).Value = 1346.5454
$ws.Cells.Item(43, 9).Value = 1381
$ws.Cells.Item(43, 10).Value = 1002
$ws.Cells.Item(43, 11).Value = 1381
$ws.Cells.Item(43, 12).Value = 1002
$ws.Cells.Item(43, 13).Value = -1312
$ws.Cells.Item(43, 14).Value = -1140
$ws.Cells.Item(58, 8).Value = 119
$ws.Cells.Item(58, 10).Value = 200
$ws.Cells.Item(58, 12).Value = 600
$ws.Cells.Item(58, 14).Value = -900
$ws.Cells.Item(70, 8).Value = 1540.0667
$ws.Cells.Item(70, 9).Value = 806.5714
$ws.Cells.Item(70, 10).Value = 2181.875
$ws.Cells.Item(70, 11).Value = 2419.7142
$ws.Cells.Item(70, 12).Value = 6545.625
$ws.Cells.Item(70, 13).Value = -2149.7142
$ws.Cells.Item(70, 14).Value = -7085.625
$ws.Cells.Item(73, 8).Value = 1540.0667
$ws.Cells.Item(73, 9).Value = 806.5714
$ws.Cells.Item(73, 10).Value = 2181.875
$ws.Cells.Item(73, 11).Value = 2419.7142
$ws.Cells.Item(73, 12).Value = 6545.625
$ws.Cells.Item(73, 13).Value = -1483.7142
$ws.Cells.Item(73, 14).Value = -8417.625
$ws.Cells.Item(86, 8).Value = 2401.6365
$ws.Cells.Item(86, 9).Value = 2274.4285
$ws.Cells.Item(86, 11).Value = 2274.4285
$ws.Cells.Item(86, 13).Value = -1151.4285
$ws.Cells.Item(89, 8).Value = 2401.6365
$ws.Cells.Item(89, 9).Value = 2274.4285
$ws.Cells.Item(89, 11).Value = 11372.1425
$ws.Cells.Item(89, 13).Value = -5756.1425
$ws.Cells.Item(135, 8).Value = 1517.2222
$ws.Cells.Item(135, 9).Value = 1276.6666
$ws.Cells.Item(135, 11).Value = 11489.9994
$ws.Cells.Item(135, 13).Value = -8954.999400000001
$ws.Cells.Item(137, 8).Value = 69277.516
$ws.Cells.Item(137, 9).Value = 121551.53
$ws.Cells.Item(137, 10).Value = 3935
$ws.Cells.Item(137, 11).Value = 364654.59
$ws.Cells.Item(137, 12).Value = 11805
$ws.Cells.Item(137, 13).Value = -362104.59
$ws.Cells.Item(137, 14).Value = -16905
$ws.Cells.Item(141, 8).Value = 2011.2858
$ws.Cells.Item(141, 9).Value = 2011.2858
$ws.Cells.Item(141, 11).Value = 6033.857400000001
$ws.Cells.Item(141, 13).Value = -853.8574000000008

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1978.5312
$ws.Cells.Item(61, 9).Value = 1506.5625
$ws.Cells.Item(61, 10).Value = 2450.5
$ws.Cells.Item(61, 11).Value = 1506.5625
$ws.Cells.Item(61, 12).Value = 2450.5
$ws.Cells.Item(61, 13).Value = -1294.5625
$ws.Cells.Item(61, 14).Value = -2874.5
$ws.Cells.Item(102, 8).Value = 3647.923
$ws.Cells.Item(102, 9).Value = 3005.389
$ws.Cells.Item(102, 11).Value = 3005.389
$ws.Cells.Item(102, 13).Value = -1383.389
$ws.Cells.Item(122, 8).Value = 8891496
$ws.Cells.Item(122, 9).Value = 11697969
$ws.Cells.Item(122, 11).Value = 35093907
$ws.Cells.Item(122, 13).Value = -35091457
$ws.Cells.Item(132, 8).Value = 1518.8605
$ws.Cells.Item(132, 9).Value = 923.3333
$ws.Cells.Item(132, 10).Value = 3484.1
$ws.Cells.Item(132, 11).Value = 2769.9999
$ws.Cells.Item(132, 12).Value = 10452.3
$ws.Cells.Item(132, 13).Value = -239.9998999999998
$ws.Cells.Item(132, 14).Value = -15512.3
$ws.Cells.Item(136, 8).Value = 1978.5312
$ws.Cells.Item(136, 9).Value = 1506.5625
$ws.Cells.Item(136, 10).Value = 2450.5
$ws.Cells.Item(136, 11).Value = 4519.6875
$ws.Cells.Item(136, 12).Value = 7351.5
$ws.Cells.Item(136, 13).Value = -1969.6875
$ws.Cells.Item(136, 14).Value = -12451.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4345.4375
$ws.Cells.Item(134, 9).Value = 2306.1052
$ws.Cells.Item(134, 11).Value = 6918.3156
$ws.Cells.Item(134, 13).Value = -4383.3156

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 27126.473
$ws.Cells.Item(31, 9).Value = 1290.2858
$ws.Cells.Item(31, 10).Value = 63297.133
$ws.Cells.Item(31, 11).Value = 1290.2858
$ws.Cells.Item(31, 12).Value = 63297.133
$ws.Cells.Item(31, 13).Value = -995.2858000000001
$ws.Cells.Item(31, 14).Value = -63887.133
$ws.Cells.Item(34, 8).Value = 27126.473
$ws.Cells.Item(34, 9).Value = 1290.2858
$ws.Cells.Item(34, 10).Value = 63297.133
$ws.Cells.Item(34, 11).Value = 1290.2858
$ws.Cells.Item(34, 12).Value = 63297.133
$ws.Cells.Item(34, 13).Value = -1088.2858
$ws.Cells.Item(34, 14).Value = -63701.133
$ws.Cells.Item(50, 8).Value = 5249.95
$ws.Cells.Item(50, 10).Value = 5249.95
$ws.Cells.Item(50, 12).Value = 5249.95
$ws.Cells.Item(50, 14).Value = -6499.95
$ws.Cells.Item(134, 8).Value = 3494.8125
$ws.Cells.Item(134, 9).Value = 2470.8
$ws.Cells.Item(134, 10).Value = 5201.5
$ws.Cells.Item(134, 11).Value = 7412.400000000001
$ws.Cells.Item(134, 12).Value = 15604.5
$ws.Cells.Item(134, 13).Value = -4877.400000000001
$ws.Cells.Item(134, 14).Value = -20674.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 2981.25
$ws.Cells.Item(140, 9).Value = 2777.6
$ws.Cells.Item(140, 11).Value = 8332.799999999999
$ws.Cells.Item(140, 13).Value = -3152.799999999999
$ws.Cells.Item(141, 8).Value = 2529.4
$ws.Cells.Item(141, 10).Value = 2750
$ws.Cells.Item(141, 12).Value = 8250
$ws.Cells.Item(141, 14).Value = -18610

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7698
$ws.Cells.Item(70, 9).Value = 8598.556
$ws.Cells.Item(70, 11).Value = 8598.556
$ws.Cells.Item(70, 13).Value = -8328.556
$ws.Cells.Item(73, 8).Value = 7698
$ws.Cells.Item(73, 9).Value = 8598.556
$ws.Cells.Item(73, 11).Value = 8598.556
$ws.Cells.Item(73, 13).Value = -7662.556
$ws.Cells.Item(80, 8).Value = 13476.5
$ws.Cells.Item(80, 10).Value = 13476.5
$ws.Cells.Item(80, 12).Value = 13476.5
$ws.Cells.Item(80, 14).Value = -15472.5
$ws.Cells.Item(83, 8).Value = 13476.5
$ws.Cells.Item(83, 10).Value = 13476.5
$ws.Cells.Item(83, 12).Value = 67382.5
$ws.Cells.Item(83, 14).Value = -77366.5
$ws.Cells.Item(122, 8).Value = 361187.9
$ws.Cells.Item(122, 9).Value = 639007.5
$ws.Cells.Item(122, 10).Value = 7599.364
$ws.Cells.Item(122, 11).Value = 1917022.5
$ws.Cells.Item(122, 12).Value = 22798.092
$ws.Cells.Item(122, 13).Value = -1914572.5
$ws.Cells.Item(122, 14).Value = -27698.092
$ws.Cells.Item(132, 8).Value = 3369.3242
$ws.Cells.Item(132, 9).Value = 2684.4583
$ws.Cells.Item(132, 11).Value = 8053.374899999999
$ws.Cells.Item(132, 13).Value = -5523.374899999999

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 2555.5
$ws.Cells.Item(68, 9).Value = 1838.8
$ws.Cells.Item(68, 10).Value = 3750
$ws.Cells.Item(68, 11).Value = 1838.8
$ws.Cells.Item(68, 12).Value = 3750
$ws.Cells.Item(68, 13).Value = -1089.8
$ws.Cells.Item(68, 14).Value = -5248
$ws.Cells.Item(71, 8).Value = 2555.5
$ws.Cells.Item(71, 9).Value = 1838.8
$ws.Cells.Item(71, 10).Value = 3750
$ws.Cells.Item(71, 11).Value = 9194
$ws.Cells.Item(71, 12).Value = 18750
$ws.Cells.Item(71, 13).Value = -5450
$ws.Cells.Item(71, 14).Value = -26238
$ws.Cells.Item(82, 8).Value = 1018.2917
$ws.Cells.Item(82, 9).Value = 836.9
$ws.Cells.Item(82, 10).Value = 1147.8572
$ws.Cells.Item(82, 11).Value = 836.9
$ws.Cells.Item(82, 12).Value = 1147.8572
$ws.Cells.Item(82, 13).Value = -475.9
$ws.Cells.Item(82, 14).Value = -1869.8572
$ws.Cells.Item(85, 8).Value = 1018.2917
$ws.Cells.Item(85, 9).Value = 836.9
$ws.Cells.Item(85, 10).Value = 1147.8572
$ws.Cells.Item(85, 11).Value = 836.9
$ws.Cells.Item(85, 12).Value = 1147.8572
$ws.Cells.Item(85, 13).Value = 411.1
$ws.Cells.Item(85, 14).Value = -3643.8572
$ws.Cells.Item(132, 8).Value = 6357.1875
$ws.Cells.Item(132, 9).Value = 7596.0586
$ws.Cells.Item(132, 10).Value = 4953.1333
$ws.Cells.Item(132, 11).Value = 22788.1758
$ws.Cells.Item(132, 12).Value = 14859.3999
$ws.Cells.Item(132, 13).Value = -20258.1758
$ws.Cells.Item(132, 14).Value = -19919.3999

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 9609.454
$ws.Cells.Item(62, 10).Value = 9609.454
$ws.Cells.Item(62, 12).Value = 9609.454
$ws.Cells.Item(62, 14).Value = -10857.454
$ws.Cells.Item(65, 8).Value = 9609.454
$ws.Cells.Item(65, 10).Value = 9609.454
$ws.Cells.Item(65, 12).Value = 48047.27
$ws.Cells.Item(65, 14).Value = -54287.27
